# Applies weekly Fruta/Hortaliza data refresh: row 2-26 field values for
# D (Fecha), L (Calidad), M (Volumen), N/O/P (precios), Q (unidad), S ($/Kg), T (Kg/unidad)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 44505
$ws.Range("L2").Value = 'Segunda'
$ws.Range("M2").Value = 300

# Row 3
$ws.Range("L3").Value = 'Primera'
$ws.Range("M3").Value = 200
$ws.Range("N3").Value = 19000
$ws.Range("O3").Value = 20000
$ws.Range("P3").Value = 19500
$ws.Range("S3").Value = 2438

# Row 4
$ws.Range("D4").Value = 44530
$ws.Range("L4").Value = 'Segunda'
$ws.Range("M4").Value = 100
$ws.Range("N4").Value = 16000
$ws.Range("O4").Value = 16000
$ws.Range("P4").Value = 16000
$ws.Range("Q4").Value = '$/bandeja 8 kilos'
$ws.Range("S4").Value = 2000
$ws.Range("T4").Value = 8

# Row 5
$ws.Range("D5").Value = 44533
$ws.Range("M5").Value = 300
$ws.Range("N5").Value = 18000
$ws.Range("O5").Value = 19000
$ws.Range("P5").Value = 18500
$ws.Range("S5").Value = 2312

# Row 6
$ws.Range("D6").Value = 44533
$ws.Range("M6").Value = 100
$ws.Range("N6").Value = 16000
$ws.Range("O6").Value = 16000
$ws.Range("P6").Value = 16000
$ws.Range("S6").Value = 2000

# Row 7
$ws.Range("D7").Value = 44498
$ws.Range("M7").Value = 300
$ws.Range("N7").Value = 19000
$ws.Range("O7").Value = 20000
$ws.Range("P7").Value = 19500
$ws.Range("S7").Value = 2438

# Row 8
$ws.Range("D8").Value = 44890
$ws.Range("M8").Value = 200
$ws.Range("N8").Value = 22000
$ws.Range("O8").Value = 22500
$ws.Range("P8").Value = 22250
$ws.Range("S8").Value = 2781

# Row 9
$ws.Range("D9").Value = 44873
$ws.Range("L9").Value = 'Primera'
$ws.Range("M9").Value = 300
$ws.Range("N9").Value = 22000
$ws.Range("O9").Value = 22500
$ws.Range("P9").Value = 22250
$ws.Range("S9").Value = 2781

# Row 10
$ws.Range("D10").Value = 44526
$ws.Range("N10").Value = 21000
$ws.Range("O10").Value = 21000
$ws.Range("P10").Value = 21000
$ws.Range("S10").Value = 2625

# Row 11
$ws.Range("D11").Value = 44491
$ws.Range("N11").Value = 18000
$ws.Range("O11").Value = 19000
$ws.Range("P11").Value = 18500
$ws.Range("S11").Value = 2312

# Row 12
$ws.Range("D12").Value = 44876
$ws.Range("L12").Value = 'Primera'
$ws.Range("N12").Value = 22000
$ws.Range("O12").Value = 22500
$ws.Range("P12").Value = 22250
$ws.Range("S12").Value = 2781

# Row 13
$ws.Range("D13").Value = 44519
$ws.Range("L13").Value = 'Primera'
$ws.Range("M13").Value = 400
$ws.Range("N13").Value = 21000
$ws.Range("O13").Value = 22000
$ws.Range("P13").Value = 21500
$ws.Range("S13").Value = 2688

# Row 14
$ws.Range("D14").Value = 44519
$ws.Range("O14").Value = 18000
$ws.Range("P14").Value = 18000
$ws.Range("S14").Value = 2250

# Row 15
$ws.Range("D15").Value = 44162
$ws.Range("M15").Value = 200

# Row 16
$ws.Range("D16").Value = 44512
$ws.Range("L16").Value = 'Segunda'
$ws.Range("N16").Value = 19000
$ws.Range("O16").Value = 20000
$ws.Range("P16").Value = 19500
$ws.Range("S16").Value = 2438

# Row 17
$ws.Range("D17").Value = 44159
$ws.Range("L17").Value = 'Primera'
$ws.Range("N17").Value = 2000
$ws.Range("O17").Value = 2100
$ws.Range("P17").Value = 2050
$ws.Range("Q17").Value = '$/kilo (en caja de 14 kilos)'
$ws.Range("S17").Value = 2050
$ws.Range("T17").Value = 1

# Row 18
$ws.Range("D18").Value = 44523
$ws.Range("L18").Value = 'Primera'
$ws.Range("M18").Value = 400
$ws.Range("N18").Value = 21000
$ws.Range("O18").Value = 22000
$ws.Range("P18").Value = 21500
$ws.Range("S18").Value = 2688

# Row 19
$ws.Range("D19").Value = 44523
$ws.Range("L19").Value = 'Segunda'
$ws.Range("M19").Value = 100
$ws.Range("N19").Value = 18000
$ws.Range("O19").Value = 18000
$ws.Range("P19").Value = 18000
$ws.Range("S19").Value = 2250

# Row 20
$ws.Range("D20").Value = 44488
$ws.Range("L20").Value = 'Segunda'
$ws.Range("M20").Value = 160
$ws.Range("N20").Value = 17000
$ws.Range("O20").Value = 18000
$ws.Range("P20").Value = 17500
$ws.Range("S20").Value = 2188

# Row 21
$ws.Range("D21").Value = 44509
$ws.Range("L21").Value = 'Segunda'
$ws.Range("N21").Value = 19000
$ws.Range("O21").Value = 20000
$ws.Range("P21").Value = 19500
$ws.Range("S21").Value = 2438

# Row 22
$ws.Range("D22").Value = 44880
$ws.Range("M22").Value = 300

# Row 23
$ws.Range("D23").Value = 44895
$ws.Range("M23").Value = 200
$ws.Range("N23").Value = 22000
$ws.Range("O23").Value = 22500
$ws.Range("P23").Value = 22250
$ws.Range("S23").Value = 2781

# Row 24
$ws.Range("D24").Value = 44894
$ws.Range("L24").Value = 'Primera'
$ws.Range("M24").Value = 200
$ws.Range("N24").Value = 22000
$ws.Range("O24").Value = 22500
$ws.Range("P24").Value = 22250
$ws.Range("S24").Value = 2781

# Row 25
$ws.Range("D25").Value = 44516
$ws.Range("L25").Value = 'Segunda'
$ws.Range("M25").Value = 200
$ws.Range("N25").Value = 18000
$ws.Range("O25").Value = 19000
$ws.Range("P25").Value = 18500
$ws.Range("S25").Value = 2312

# Row 26
$ws.Range("D26").Value = 44495
$ws.Range("M26").Value = 270
$ws.Range("N26").Value = 19000
$ws.Range("O26").Value = 20000
$ws.Range("P26").Value = 19556
$ws.Range("S26").Value = 2444

